$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix of the energy parameters (column A values)
$ws.Range("A3").Value = 1.42
$ws.Range("A5").Value = 1.85
$ws.Range("A6").Value = 1.98
$ws.Range("A7").Value = 2.7

# New cell with a single space string value
$ws.Range("E14").Value = " "

# Page setup (adds <pageSetup .../> to the sheet XML, matching the
# author's printer/page-setup touch)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# View state: zoom level and active selection as left by the author
$excel.ActiveWindow.Zoom = 161
[void]$ws.Range("F12").Select()
